# Apply crypto price/volume updates from GitHub Actions scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '23.306.97'
$ws.Cells.Item(2, 5).Value = '  -0.46%  '
$ws.Cells.Item(3, 4).Value = '1.624.12'
$ws.Cells.Item(3, 5).Value = '  -0.92%  '
$ws.Cells.Item(4, 5).Value = '  +0.04%  '
$ws.Cells.Item(5, 5).Value = '  -0.02%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '302.69'
$ws.Cells.Item(6, 5).Value = '  -0.71%  '
$ws.Cells.Item(7, 5).Value = '  +0.51%  '
$ws.Cells.Item(8, 5).Value = '  -0.29%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '51.24'
$ws.Cells.Item(9, 5).Value = '  -1.79%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.08129'
$ws.Cells.Item(10, 5).Value = '  +0.12%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '1.218'
$ws.Cells.Item(11, 5).Value = '  -2.48%  '
$ws.Cells.Item(12, 5).Value = '  +0.04%  '
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '22.23'
$ws.Cells.Item(13, 5).Value = '  -2.74%  '
$ws.Cells.Item(14, 5).Value = '  -2.06%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '0.00001236'
$ws.Cells.Item(15, 5).Value = '  -2.76%  '
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '7.265'
$ws.Cells.Item(16, 5).Value = '  -0.14%  '
$ws.Cells.Item(17, 4).Value = '1.614.59'
$ws.Cells.Item(17, 5).Value = '  -0.88%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '94.04'
$ws.Cells.Item(18, 5).Value = '  -0.38%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '0.06934'
$ws.Cells.Item(19, 5).Value = '  +0.85%  '
$ws.Cells.Item(20, 2).Value = 'Uniswap'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '6.551'
$ws.Cells.Item(20, 5).Value = '  +0.56%  '
$ws.Cells.Item(21, 2).Value = 'Avalanche'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '17.46'
$ws.Cells.Item(21, 5).Value = '  -3.68%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '12.45'
$ws.Cells.Item(23, 5).Value = '  -2.31%  '
$ws.Cells.Item(24, 4).Value = '23.299.00'
$ws.Cells.Item(24, 5).Value = '  -0.56%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '2.493'
$ws.Cells.Item(25, 5).Value = '  +3.31%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '3.080'
$ws.Cells.Item(26, 5).Value = '  +1.98%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '21.09'
$ws.Cells.Item(27, 5).Value = '  -0.54%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '150.45'
$ws.Cells.Item(28, 5).Value = '  -0.98%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '5.287'
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '132.79'
$ws.Cells.Item(30, 5).Value = '  -1.55%  '
$ws.Cells.Item(31, 4).Value = '1.796.50'
$ws.Cells.Item(31, 5).Value = '  -0.68%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '6.707'
$ws.Cells.Item(32, 5).Value = '  -0.90%  '
$ws.Cells.Item(33, 2).Value = 'ImmutableX'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '1.059'
$ws.Cells.Item(33, 5).Value = '  +11.23%  '
$ws.Cells.Item(34, 2).Value = 'WEMIXTOKEN'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '2.086'
$ws.Cells.Item(34, 5).Value = '  -8.99%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '11.21'
$ws.Cells.Item(35, 5).Value = '  +8.62%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.02739'
$ws.Cells.Item(36, 5).Value = '  -3.35%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.08760'
$ws.Cells.Item(37, 5).Value = '  -0.23%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.2473'
$ws.Cells.Item(38, 5).Value = '  -1.99%  '
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.07085'
$ws.Cells.Item(39, 5).Value = '  -1.68%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '5.972'
$ws.Cells.Item(40, 5).Value = '  -1.54%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.6950'
$ws.Cells.Item(41, 5).Value = '  -1.28%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '1.327'
$ws.Cells.Item(42, 5).Value = '  -3.56%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '15.80'
$ws.Cells.Item(43, 5).Value = '  -0.99%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '12.00'
$ws.Cells.Item(44, 5).Value = '  -3.52%  '
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.6433'
$ws.Cells.Item(45, 5).Value = '  -0.97%  '
$ws.Cells.Item(46, 5).Value = '  -0.01%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '2.262'
$ws.Cells.Item(47, 5).Value = '  -2.87%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '3.953'
$ws.Cells.Item(48, 5).Value = '  -1.34%  '
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.07959'
$ws.Cells.Item(49, 5).Value = '  -0.05%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '125.76'
$ws.Cells.Item(50, 5).Value = '  -2.22%  '
$ws.Cells.Item(51, 5).Value = '  -1.23%  '
